$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "GF #" (column J) values for each closing row.
# Numeric GF numbers are entered as numbers; rows that did not close
# get the literal text "did not close".
$ws.Range("J2").Value = 20191011
$ws.Range("J3").Value = 20191062
$ws.Range("J4").Value = "did not close"
$ws.Range("J5").Value = "did not close"
$ws.Range("J6").Value = "did not close"
$ws.Range("J7").Value = "did not close"
$ws.Range("J8").Value = 20190096
$ws.Range("J9").Value = "did not close"
$ws.Range("J10").Value = 20182625
$ws.Range("J11").Value = "did not close"

# Mirror the author's final selection: the whole new column of data,
# anchored at the last entered cell.
$ws.Range("J2:J11").Select()
